$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: assign a value while forcing the cell to stay Text-typed
# (the source data stores every Coin/Link/Price/Volume cell as an inline
# string, even numeric-looking prices like "227.93" - without this, Excel
# would silently reinterpret a clean numeric string as a Number cell).
function Set-TextValue($range, $value) {
    $r = $ws.Range($range)
    $origStyle = $r.Style
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.Style = $origStyle
}

$updates = @(
    @{ Row = 2; D = '37.796.23'; E = '  +0.38%  ' },
    @{ Row = 3; D = '2.025.67'; E = '  -0.65%  ' },
    @{ Row = 4; E = '  +0.00%  ' },
    @{ Row = 5; D = '227.06'; E = '  -1.19%  ' },
    @{ Row = 6; D = '0.611'; E = '  -0.22%  ' },
    @{ Row = 7; D = '60.01'; E = '  +6.59%  ' },
    @{ Row = 8; E = '  +0.00%  ' },
    @{ Row = 9; D = '0.384'; E = '  +0.34%  ' },
    @{ Row = 10; D = '0.0805'; E = '  +0.32%  ' },
    @{ Row = 11; E = '  +0.68%  ' },
    @{ Row = 12; D = '14.55'; E = '  +1.01%  ' },
    @{ Row = 13; D = '2.326.62'; E = '  -0.72%  ' },
    @{ Row = 14; D = '21.20'; E = '  +4.32%  ' },
    @{ Row = 15; D = '0.751'; E = '  +1.25%  ' },
    @{ Row = 16; D = '5.23'; E = '  +0.01%  ' },
    @{ Row = 17; D = '2.037.17'; E = '  -0.16%  ' },
    @{ Row = 18; D = '37.786.06'; E = '  +0.53%  ' },
    @{ Row = 19; D = '6.03'; E = '  -2.23%  ' },
    @{ Row = 20; D = '69.51'; E = '  +0.48%  ' },
    @{ Row = 21; D = '0.0₃0821'; E = '  -0.44%  ' },
    @{ Row = 22; D = '224.41'; E = '  +0.56%  ' },
    @{ Row = 23; E = '  +0.12%  ' },
    @{ Row = 24; D = '2.42'; E = '  -0.97%  ' },
    @{ Row = 25; E = '  -2.25%  ' },
    @{ Row = 26; D = '165.34'; E = '  +0.05%  ' },
    @{ Row = 27; D = '9.18'; E = '  +0.08%  ' },
    @{ Row = 28; D = '0.130'; E = '  -1.69%  ' },
    @{ Row = 29; D = '18.83'; E = '  -0.68%  ' },
    @{ Row = 30; D = '1.28'; E = '  -3.93%  ' },
    @{ Row = 31; E = '  +1.55%  ' },
    @{ Row = 32; D = '4.43'; E = '  -1.05%  ' },
    @{ Row = 33; D = '2.04'; E = '  +1.10%  ' },
    @{ Row = 34; D = '4.49'; E = '  -0.09%  ' },
    @{ Row = 35; D = '0.0601'; E = '  -1.20%  ' },
    @{ Row = 36; D = '6.29'; E = '  +6.63%  ' },
    @{ Row = 37; D = '2.25'; E = '  -3.85%  ' },
    @{ Row = 38; E = '  +0.29%  ' },
    @{ Row = 39; E = '  +0.02%  ' },
    @{ Row = 40; D = '1.532.90'; E = '  +3.84%  ' },
    @{ Row = 41; D = '0.0217'; E = '  +0.43%  ' },
    @{ Row = 42; D = '96.21'; E = '  +0.30%  ' },
    @{ Row = 43; D = '16.67'; E = '  +0.68%  ' },
    @{ Row = 44; D = '0.0917'; E = '  -3.28%  ' },
    @{ Row = 45; D = '2.79'; E = '  -2.11%  ' },
    @{ Row = 46; D = '1.10'; E = '  -0.74%  ' },
    @{ Row = 47; D = '3.96'; E = '  -3.87%  ' },
    @{ Row = 48; B = 'ARBITRUM'; C = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'; D = '1.00'; E = '  -1.01%  ' },
    @{ Row = 49; B = 'MXToken'; C = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'; D = '2.96'; E = '  +0.75%  ' },
    @{ Row = 50; D = '7.07'; E = '  -0.67%  ' },
    @{ Row = 51; D = '2.215.25'; E = '  -0.68%  ' }
)

foreach ($u in $updates) {
    if ($u.ContainsKey("B")) { $ws.Range("B" + $u.Row).Value = $u.B }
    if ($u.ContainsKey("C")) { $ws.Range("C" + $u.Row).Value = $u.C }
    if ($u.ContainsKey("D")) { Set-TextValue ("D" + $u.Row) $u.D }
    if ($u.ContainsKey("E")) { $ws.Range("E" + $u.Row).Value = $u.E }
}
